# PST 4.0.0 and UM 1.0.0 Release
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates ---

# Row 4: "Version" value becomes a plain number (3.1) instead of the text "3.1.1"
$ws.Range("C4").Value = 3.1

# Row 4: Sim. Time [s]
$ws.Range("S4").Value = 740.35

# Row 5: Sim. Time [s] (U5's shared formula $S$4/S5 recalculates automatically)
$ws.Range("S5").Value = 284.37

# Row 6: Step size / solutions-per-step / total-steps / sim-time updates
$ws.Range("G6").Value = 0.000136
$ws.Range("I6").Value = 0.0138
$ws.Range("K6").Value = 17353
$ws.Range("O6").Value = 96
$ws.Range("Q6").Value = 27243
$ws.Range("S6").Value = 125.53

# --- View state: move the selection/scroll position ---
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S9").Select() | Out-Null
